# Swap the test-step content between TC3 (row 26) and TC4 (row 33):
#  - TC3's step/expected-result becomes the former TC4 "atribuir/desatribuir" content
#  - TC4's step/expected-result becomes the former TC3 "realizar o empenho" content
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tc3Step   = $ws.Range("B26").Value2
$tc3Result = $ws.Range("D26").Value2
$tc4Step   = $ws.Range("B33").Value2
$tc4Result = $ws.Range("D33").Value2

$ws.Range("B26").Value2 = $tc4Step
$ws.Range("D26").Value2 = $tc4Result
$ws.Range("B33").Value2 = $tc3Step
$ws.Range("D33").Value2 = $tc3Result
